$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.850.82"
$ws.Range("E2").Value = "  -5.73%  "
$ws.Range("D3").Value = "3.280.74"
$ws.Range("E3").Value = "  -7.19%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'519.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.59%  "
$ws.Range("D6").Value = "'172.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -15.14%  "
$ws.Range("D7").Value = "'0.600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("D8").Value = "3.277.13"
$ws.Range("E8").Value = "  -7.05%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'0.602"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.78%  "
$ws.Range("D11").Value = "'55.79"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.76%  "
$ws.Range("E12").Value = "  -8.32%  "
$ws.Range("D13").Value = "'0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.01%  "
$ws.Range("D14").Value = "'8.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.81%  "
$ws.Range("D15").Value = "3.799.31"
$ws.Range("E15").Value = "  -7.28%  "
$ws.Range("D16").Value = "3.269.91"
$ws.Range("E16").Value = "  -7.56%  "
$ws.Range("E17").Value = "  -6.97%  "
$ws.Range("D18").Value = "63.729.25"
$ws.Range("E18").Value = "  -5.63%  "
$ws.Range("D19").Value = "'17.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.91%  "
$ws.Range("D20").Value = "'10.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.98%  "
$ws.Range("D21").Value = "'0.947"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.30%  "
$ws.Range("D22").Value = "'370.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.08%  "
$ws.Range("D23").Value = "'3.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.76%  "
$ws.Range("D24").Value = "'79.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.34%  "
$ws.Range("D25").Value = "'10.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.82%  "
$ws.Range("D26").Value = "'3.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("D28").Value = "'2.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.39%  "
$ws.Range("D29").Value = "'11.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.96%  "
$ws.Range("D30").Value = "'8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.03%  "
$ws.Range("D31").Value = "'28.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.75%  "
$ws.Range("D32").Value = "'636.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -12.51%  "
$ws.Range("D33").Value = "'6.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.60%  "
$ws.Range("D34").Value = "'11.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'58.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -8.63%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.104"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.26%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'36.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.80%  "
$ws.Range("E39").Value = "  -4.92%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "2.900.34"
$ws.Range("E42").Value = "  -5.86%  "
$ws.Range("D43").Value = "'0.121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.72%  "
$ws.Range("D44").Value = "'2.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.17%  "
$ws.Range("E45").Value = "  -14.83%  "
$ws.Range("E46").Value = "  -4.81%  "
$ws.Range("D47").Value = "'0.0392"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.28%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.124"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'2.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.23%  "
$ws.Range("E50").Value = "  +3.08%  "
$ws.Range("D51").Value = "'134.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.24%  "
